$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet had a blank separator row at 22 (just under the last data row,
# 21) followed by the summary rows (sum [min] / sum [h] / sum [working
# weeks]) at 23-25. A new data row is being inserted at 22, which pushes
# the blank separator + summary rows down by one (23/24/25 -> 24/25/26)
# while keeping their formulas/labels intact (Excel adjusts the SUM/ratio
# formulas' relative references automatically on insert).
$ws.Rows("22:22").Insert()

# New data row: 2014-02-24, 10:30 -> 12:30.
$ws.Range("A22").Value = 2014
$ws.Range("B22").Value = 2
$ws.Range("C22").Value = 24
$ws.Range("D22").Value = 0.4375
$ws.Range("E22").Value = 0.52083333333333337

# Match the time/number formatting used by the rest of the data rows.
$ws.Range("D22").NumberFormat = $ws.Range("D21").NumberFormat
$ws.Range("E22").NumberFormat = $ws.Range("E21").NumberFormat
$ws.Range("F22").NumberFormat = $ws.Range("F21").NumberFormat
$ws.Range("G22").NumberFormat = $ws.Range("G21").NumberFormat

# Continue the "time spent" formulas into the new row.
$ws.Range("F22").Formula = "=(E22-D22)*24*60"
$ws.Range("G22").Formula = "=F22/60"

# The diff moves the active selection from E22 to F22.
$ws.Range("F22").Select()

$wb.Save()
